$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new "Save" column (H1), copying the formatting of the
# existing header cell G1 ("sum") so it matches the other header cells.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add the new column's data values
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
